$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated strikeout (K) values for rows 2-35 (column G), replacing the
# previous "Strike#" derived values with the regenerated K values.
$sVals = @(5,5,2,2,1,5,5,5,2,4,4,1,5,2,3,3,5,1,5,0,6,1,1,4,5,2,5,2,2,7,7,6,1,2)

$startRow = 2
for ($i = 0; $i -lt $sVals.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $sVals[$i]
}
